$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2659.4119
$ws.Range("J17").Value = 2659.4119
$ws.Range("L17").Value = 7978.2357
$ws.Range("N17").Value = -8314.235700000001

$ws.Range("H76").Value = 4825.3335
$ws.Range("I76").Value = 5211.1113
$ws.Range("J76").Value = 3668
$ws.Range("K76").Value = 5211.1113
$ws.Range("L76").Value = 3668
$ws.Range("M76").Value = -4896.1113
$ws.Range("N76").Value = -4298

$ws.Range("H79").Value = 4825.3335
$ws.Range("I79").Value = 5211.1113
$ws.Range("J79").Value = 3668
$ws.Range("K79").Value = 5211.1113
$ws.Range("L79").Value = 3668
$ws.Range("M79").Value = -4119.1113
$ws.Range("N79").Value = -5852

$ws.Range("H100").Value = 933.7059
$ws.Range("I100").Value = 612.3570999999999
$ws.Range("J100").Value = 2433.3333
$ws.Range("K100").Value = 612.3570999999999
$ws.Range("L100").Value = 2433.3333
$ws.Range("M100").Value = -71.35709999999995
$ws.Range("N100").Value = -3515.3333

$ws.Range("H132").Value = 11119810
$ws.Range("I132").Value = 18527754
$ws.Range("J132").Value = 7894.0835
$ws.Range("K132").Value = 55583262
$ws.Range("L132").Value = 23682.2505
$ws.Range("M132").Value = -55580732
$ws.Range("N132").Value = -28742.2505

$ws.Range("H137").Value = 1439.9231
$ws.Range("J137").Value = 1808.5
$ws.Range("L137").Value = 5425.5
$ws.Range("N137").Value = -10525.5

$ws.Range("H138").Value = 552017.5600000001
$ws.Range("I138").Value = 1597.75
$ws.Range("J138").Value = 650600.2
$ws.Range("K138").Value = 4793.25
$ws.Range("L138").Value = 1951800.6
$ws.Range("M138").Value = 346.75
$ws.Range("N138").Value = -1962080.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 66667910
$ws.Range("I61").Value = 76924130
$ws.Range("K61").Value = 76924130
$ws.Range("M61").Value = -76923918

$ws.Range("H74").Value = 1807.1818
$ws.Range("I74").Value = 1343.8667
$ws.Range("J74").Value = 2800
$ws.Range("K74").Value = 1343.8667
$ws.Range("L74").Value = 2800
$ws.Range("M74").Value = -469.8667
$ws.Range("N74").Value = -4548

$ws.Range("H77").Value = 1807.1818
$ws.Range("I77").Value = 1343.8667
$ws.Range("J77").Value = 2800
$ws.Range("K77").Value = 6719.333500000001
$ws.Range("L77").Value = 14000
$ws.Range("M77").Value = -2351.333500000001
$ws.Range("N77").Value = -22736

$ws.Range("H136").Value = 66667910
$ws.Range("I136").Value = 76924130
$ws.Range("K136").Value = 230772390
$ws.Range("M136").Value = -230769840

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 430.46155
$ws.Range("J80").Value = 469.6
$ws.Range("L80").Value = 469.6
$ws.Range("N80").Value = -2465.6

$ws.Range("H83").Value = 430.46155
$ws.Range("J83").Value = 469.6
$ws.Range("L83").Value = 2348
$ws.Range("N83").Value = -12332

$ws.Range("H116").Value = 59999.5
$ws.Range("J116").Value = 59999.5
$ws.Range("L116").Value = 59999.5
$ws.Range("N116").Value = -69177.5

$ws.Range("H134").Value = 1453.8667
$ws.Range("I134").Value = 1181.3
$ws.Range("K134").Value = 3543.9
$ws.Range("M134").Value = -1008.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1380.9756
$ws.Range("I31").Value = 1357.9487
$ws.Range("J31").Value = 1830
$ws.Range("K31").Value = 1357.9487
$ws.Range("L31").Value = 1830
$ws.Range("M31").Value = -1062.9487
$ws.Range("N31").Value = -2420

$ws.Range("H34").Value = 1380.9756
$ws.Range("I34").Value = 1357.9487
$ws.Range("J34").Value = 1830
$ws.Range("K34").Value = 1357.9487
$ws.Range("L34").Value = 1830
$ws.Range("M34").Value = -1155.9487
$ws.Range("N34").Value = -2234

$ws.Range("H58").Value = 1588.0286
$ws.Range("I58").Value = 1264
$ws.Range("J58").Value = 2020.0667
$ws.Range("K58").Value = 1264
$ws.Range("L58").Value = 2020.0667
$ws.Range("M58").Value = -1061
$ws.Range("N58").Value = -2426.0667

$ws.Range("H136").Value = 1588.0286
$ws.Range("I136").Value = 1264
$ws.Range("J136").Value = 2020.0667
$ws.Range("K136").Value = 3792
$ws.Range("L136").Value = 6060.2001
$ws.Range("M136").Value = -1242
$ws.Range("N136").Value = -11160.2001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H115").Value = 5072.923
$ws.Range("I115").Value = 3414
$ws.Range("J115").Value = 5374.5454
$ws.Range("K115").Value = 10242
$ws.Range("L115").Value = 16123.6362
$ws.Range("M115").Value = -9067
$ws.Range("N115").Value = -18473.6362

$ws.Range("H131").Value = 15625887
$ws.Range("J131").Value = 1059.4445
$ws.Range("L131").Value = 3178.3335
$ws.Range("N131").Value = -13258.3335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 604.6667
$ws.Range("I97").Value = 590.2308
$ws.Range("K97").Value = 590.2308
$ws.Range("M97").Value = -94.23080000000004

$ws.Range("H132").Value = 2797.879
$ws.Range("I132").Value = 2562.1428
$ws.Range("J132").Value = 3210.4167
$ws.Range("K132").Value = 7686.428400000001
$ws.Range("L132").Value = 9631.250100000001
$ws.Range("M132").Value = -5156.428400000001
$ws.Range("N132").Value = -14691.2501

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 21000
$ws.Range("J63").Value = 21000
$ws.Range("L63").Value = 21000
$ws.Range("N63").Value = -22498

$ws.Range("H66").Value = 21000
$ws.Range("J66").Value = 21000
$ws.Range("L66").Value = 63000
$ws.Range("N66").Value = -70488

$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws.Range("H100").Value = 1188.909
$ws.Range("I100").Value = 1039.7142
$ws.Range("K100").Value = 1039.7142
$ws.Range("M100").Value = -498.7141999999999

$ws.Range("H122").Value = 2610.6667
$ws.Range("I122").Value = 2664.75
$ws.Range("J122").Value = 2502.5
$ws.Range("K122").Value = 7994.25
$ws.Range("L122").Value = 7507.5
$ws.Range("M122").Value = -5544.25
$ws.Range("N122").Value = -12407.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 16450
$ws.Range("J64").Value = 16450
$ws.Range("L64").Value = 16450
$ws.Range("N64").Value = -16946

$ws.Range("H67").Value = 16450
$ws.Range("J67").Value = 16450
$ws.Range("L67").Value = 16450
$ws.Range("N67").Value = -18166

$ws.Range("H121").Value = 30000
$ws.Range("J121").Value = 30000
$ws.Range("L121").Value = 30000
$ws.Range("N121").Value = -33494

$ws.Range("H123").Value = 31381.666
$ws.Range("J123").Value = 31381.666
$ws.Range("L123").Value = 31381.666
$ws.Range("N123").Value = -41181.666

$ws.Range("H126").Value = 4618.875
$ws.Range("I126").Value = 2908.5833
$ws.Range("K126").Value = 8725.749899999999
$ws.Range("M126").Value = -6255.749899999999
